# "Read through start of results." — fill in the two antibody-table cells
# that were left blank (Flt3's Clone/Supplier) and add the missing
# Streptavidin Clone placeholder ("-") on slide 3 of the "Antibody clones"
# deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The slide has a single shape: the antibody/clone/supplier table.
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

# Row 9 = "Flt3" — Clone column (2) and Supplier column (3) were empty.
$tbl.Cell(9, 2).Shape.TextFrame.TextRange.Text = "A2F10"
$tbl.Cell(9, 3).Shape.TextFrame.TextRange.Text = "eBioscience"

# Row 18 = "Streptavidin" — Clone column (2) was empty, add placeholder "-".
$tbl.Cell(18, 2).Shape.TextFrame.TextRange.Text = "-"
